# Added steel unit data from IEAGHG2013
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New unit process rows pulled from the IEAGHG2013 steel dataset.
$newRows = @(
    @{ Row = 12; ID = "IEAGHG_BOF";     Name = "basic oxygen furnace"; Product = "raw liquid steel"; Sheet = "BOF Steelmaking" },
    @{ Row = 13; ID = "IEAGHG_ladle";   Name = "ladle metallurgy";     Product = "liquid steel";      Sheet = "Ladle Metallurgy" },
    @{ Row = 14; ID = "IEAGHG_forming"; Name = "rolling & casting";    Product = "hot rolled coil";   Sheet = "Forming" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A${row}").Value = $r.ID
    $ws.Range("B${row}").Value = $r.Name
    $ws.Range("C${row}").Value = $r.Product
    $ws.Range("D${row}").Value = "outflow"
    $ws.Range("E${row}").Value = "data/steel/SteelUnits_Variables.xlsx"
    $ws.Range("F${row}").Value = $r.Sheet
    $ws.Range("G${row}").Value = "data/steel/SteelUnits_Relationships.xlsx"
    $ws.Range("H${row}").Value = $r.Sheet

    # Columns A-E and G carry the workbook's "Text" number format (style index 1);
    # F and H are left as general/default format, matching the existing rows.
    $ws.Range("A${row}:E${row}").NumberFormat = "@"
    $ws.Range("G${row}").NumberFormat = "@"
}

# Restore the view state captured in the saved workbook.
$ws.Activate()
$excel.ActiveWindow.Zoom = 150
[void]$ws.Range("C9").Select()
